# Apply the edits described by the commit diff to the "Ref_incandescent phase out" / SK workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the Branch value in A3 (shared string update):
#    "CIMS.CAN.SK.Residential.Buildings.Floorspace.Lighting"
#    -> "CIMS.CAN.SK.Residential.Dwellings.Lighting"
$ws.Range("A3").Value = "CIMS.CAN.SK.Residential.Dwellings.Lighting"

# 2. Replace the "Market share new_max" time series (M3:W3) with new, explicit
#    values. This also removes the shared formula (N3:W3 = M3) that used to
#    propagate the value across the row, since every cell now carries its own
#    literal number.
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0001
$ws.Range("S3").Value = 0.0001
$ws.Range("T3").Value = 0.0001
$ws.Range("U3").Value = 0.0001
$ws.Range("V3").Value = 0.0001
$ws.Range("W3").Value = 0.0001

# 3. Update the sheet's saved selection to match the new extent (A1:X3 -> A1:X4).
$ws.Range("A1:X4").Select()

# 4. Match the workbook's "calculate before save" behaviour (calcOnSave="0").
$excel.CalculateBeforeSave = $false

$wb.Save()
